$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the split "Testing four routers in series (1-2-3-4)" runs
#    into a single run, while preserving the <w:lastRenderedPageBreak/>
#    that lives in the paragraph's first run. A plain Find/Replace (or
#    Range.Text assignment) that touches that first run strips the
#    <w:lastRenderedPageBreak/> marker, so we rebuild the paragraph via
#    InsertXML, which round-trips raw OOXML (including the marker)
#    faithfully.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text.TrimEnd()
    if ($t -eq "Testing four routers in series (1-2-3-4)" -or
        ($t.StartsWith("Testing") -and $t.Contains("four") -and $t.Contains("routers in series (1-2-3") -and $t.Contains("4)"))) {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Exclude the trailing paragraph mark from the replaced span.
    $pr = $d.Range($r.Start, $r.End - 1)
    $mergedXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:lastRenderedPageBreak/><w:t>Testing four routers in series (1-2-3-4)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $pr.InsertXML($mergedXml)
}

# ---------------------------------------------------------------------
# 2) Append the new "four routers ... with varying cost" section after
#    the very last paragraph in the document body.
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$endPos = $lastPara.Range.End
$insertionPoint = $d.Range($endPos, $endPos)

$appendXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>
  <w:r><w:t>Testing four routers in series (1-2-3-4)</w:t></w:r>
  <w:r><w:t xml:space="preserve"> with varying cost</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>router-id 1</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>input-ports 10001</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>outputs 2001-1-2</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>router-id 2</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>input-ports 2001,2002</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>outputs 3001-2-3, 10001-1-1</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>router-id 3</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>input-ports 3001, 3002</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>outputs 2002-2-2, 4001-3-4</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>router-id 4</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>input-ports 4001</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/></w:pPr>
  <w:r><w:t>outputs 3002-3-3</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertionPoint.InsertXML($appendXml)
